# Update "want to go" (想去人数) counts in column F for the
# 展览 (Exhibitions) and 全部类型 (All types) sheets.
# Values bumped per the upstream data refresh (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 51
$wsExhibitions.Range("F3").Value = 114
$wsExhibitions.Range("F4").Value = 1571
$wsExhibitions.Range("F5").Value = 268
$wsExhibitions.Range("F7").Value = 1408
$wsExhibitions.Range("F8").Value = 10193
$wsExhibitions.Range("F10").Value = 137
$wsExhibitions.Range("F13").Value = 391
$wsExhibitions.Range("F14").Value = 7069
$wsExhibitions.Range("F15").Value = 1101
$wsExhibitions.Range("F16").Value = 663
$wsExhibitions.Range("F17").Value = 32
$wsExhibitions.Range("F19").Value = 233

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 51
$wsAllTypes.Range("F3").Value = 114
$wsAllTypes.Range("F4").Value = 1571
$wsAllTypes.Range("F5").Value = 268
$wsAllTypes.Range("F8").Value = 1408
$wsAllTypes.Range("F11").Value = 10193
$wsAllTypes.Range("F13").Value = 137
$wsAllTypes.Range("F16").Value = 391
$wsAllTypes.Range("F17").Value = 7069
$wsAllTypes.Range("F18").Value = 1101
$wsAllTypes.Range("F19").Value = 663
$wsAllTypes.Range("F20").Value = 32
$wsAllTypes.Range("F22").Value = 233
